$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the insertion point for the newly-added player "Alex Len" -
# every existing roster row from row 2 down needs to move one row lower
# to make room for it. Find the last used row first.
$lastRow = $ws.UsedRange.Rows.Count

# Shift the existing data down by one row, working from the bottom up so
# a source row is never overwritten before it has been copied.
for ($r = $lastRow; $r -ge 2; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value = $ws.Cells.Item($r, 3).Value2
}

# Write the newly added player into the now-vacant row 2.
$ws.Cells.Item(2, 1).Value = "Alex Len"
$ws.Cells.Item(2, 2).Value = "C"
$ws.Cells.Item(2, 3).Value = "Los Angeles Lakers"
